$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsSchedule.Range("B2").Value = 46040.85416666666
$wsSchedule.Range("C2").Value = 14
$wsSchedule.Range("D2").Value = 52.91999999999999
$wsSchedule.Range("E2").Value = 149.216769
$wsSchedule.Range("F2").Value = 2.819666836734695
$wsSchedule.Range("A3").Value = 46040.91666666666
$wsSchedule.Range("C3").Value = 5
$wsSchedule.Range("D3").Value = 18.9
$wsSchedule.Range("E3").Value = 520.6762177499999
$wsSchedule.Range("F3").Value = 27.54900623015872
$wsSchedule.Range("E4").Value = 20.29469325000004
$wsSchedule.Range("F4").Value = 0.5965518298059976
$wsDetailed.Range("B31").Value = -23.5
$wsDetailed.Range("B32").Value = -14.53915
$wsDetailed.Range("B33").Value = -11.93964
$wsDetailed.Range("B34").Value = -5.43826
$wsDetailed.Range("C34").Value = "historical"
$wsDetailed.Range("B35").Value = -5.41
$wsDetailed.Range("B36").Value = 8.2301
$wsDetailed.Range("B37").Value = 9.78525
$wsDetailed.Range("B38").Value = 14.2512
$wsDetailed.Range("B39").Value = 29.00544
$wsDetailed.Range("B40").Value = 30.91231
$wsDetailed.Range("B41").Value = 56.98
$wsDetailed.Range("B42").Value = 49.36694
$wsDetailed.Range("E42").Value = "ON"
$wsDetailed.Range("B44").Value = 56.98
$wsDetailed.Range("B45").Value = 56.98
$wsDetailed.Range("E45").Value = "OFF"
$wsDetailed.Range("B46").Value = 56.98
$wsDetailed.Range("B47").Value = 50.56345
$wsDetailed.Range("B49").Value = 49.15555
$wsDetailed.Range("B50").Value = 56.20781
$wsDetailed.Range("B51").Value = 56.98
$wsDetailed.Range("B52").Value = 57.06003
$wsDetailed.Range("B53").Value = 57.06003
$wsDetailed.Range("B54").Value = 56.97996
$wsDetailed.Range("B55").Value = 56.97996
$wsDetailed.Range("B56").Value = 57.06003
$wsDetailed.Range("B57").Value = 57.06003
$wsDetailed.Range("B58").Value = 57.06003
$wsDetailed.Range("B59").Value = 59.1451
$wsDetailed.Range("B60").Value = 57.99562
$wsDetailed.Range("B61").Value = 65
$wsDetailed.Range("B62").Value = 65
$wsDetailed.Range("B64").Value = 35.88
$wsDetailed.Range("B66").Value = 23.61013
$wsDetailed.Range("B67").Value = 0.7
$wsDetailed.Range("B70").Value = -6.23257
$wsDetailed.Range("B71").Value = -6.09194
$wsDetailed.Range("B72").Value = -6.17304
$wsDetailed.Range("B73").Value = -6.02722
$wsDetailed.Range("B74").Value = -6.94065
$wsDetailed.Range("B75").Value = -6.98156
$wsDetailed.Range("B76").Value = -7.47989
$wsDetailed.Range("B77").Value = -6.06786
$wsDetailed.Range("B78").Value = -5.95332
$wsDetailed.Range("B79").Value = -5.51
$wsDetailed.Range("B80").Value = -5.50985
$wsDetailed.Range("B81").Value = -0.95731
$wsDetailed.Range("B82").Value = 0.00025
$wsDetailed.Range("B83").Value = -2.54783
$wsDetailed.Range("B84").Value = 0.01029
$wsDetailed.Range("B85").Value = 5.04892
$wsDetailed.Range("B86").Value = 20.67051
$wsDetailed.Range("B87").Value = 55.45586
$wsDetailed.Range("B88").Value = 57.40696
$wsDetailed.Range("B89").Value = 73.20007
$wsDetailed.Range("B90").Value = 68.75594
$wsDetailed.Range("B91").Value = 67.64812000000001
$wsDetailed.Range("B92").Value = 64.6053
$wsDetailed.Range("B93").Value = 65
$wsDetailed.Range("B94").Value = 62.95723
$wsDetailed.Range("B95").Value = 57.98162
$wsDetailed.Range("B96").Value = 57.3
$wsDetailed.Range("B97").Value = 60.2337

Write-Host "Applied all changes"